# moved convert logic from function to class
# The "Fuel Type" value for the last fleet row (H5) is no longer populated
# inline on the sheet - clear it out, then leave the selection where the
# user's cursor ended up (H5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").ClearContents()
$ws.Range("H5").Select()
